# Add a new worksheet "2016-10-21" after the existing "2016-10-07" sheet
# and populate it with the new campus-recruitment listing rows.
#
# The new sheet is created by copying the existing sheet (so it inherits the
# same sheet-level properties: outline settings, column-width baseline, page
# margins, etc.) and then replacing its contents with the 2016-10-21 data.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Copy([System.Reflection.Missing]::Value, $sheet1)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2016-10-21"

# Wipe the copied data so we can write the new sheet's rows from scratch.
$ws.Range("A1:H51").ClearContents()

# Helper: write a value into a cell as literal text, even if it looks like
# a date (e.g. "2016-10-24"), so it lands as a shared string rather than
# being auto-converted to a numeric date serial by Excel's input parser.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Header row
$ws.Cells.Item(1, 1).Value = "number"
$ws.Cells.Item(1, 2).Value = "deadline"
$ws.Cells.Item(1, 3).Value = "location"
$ws.Cells.Item(1, 4).Value = "department"
$ws.Cells.Item(1, 5).Value = "title"
$ws.Cells.Item(1, 6).Value = "CM"
$ws.Cells.Item(1, 7).Value = "company"
$ws.Cells.Item(1, 8).Value = "link"

# Row 2 - CMB Huizhou branch posting
$ws.Cells.Item(2, 1).Value = "Not available"
Set-TextValue $ws.Cells.Item(2, 2) "2016-10-24"
$ws.Cells.Item(2, 3).Value = "惠州"
$ws.Cells.Item(2, 4).Value = "深圳分行"
$ws.Cells.Item(2, 5).Value = "招商银行惠州分行2017年秋季校园招聘"
$ws.Cells.Item(2, 6).Value = "C"
$ws.Cells.Item(2, 7).Value = "cmbchina"
$ws.Cells.Item(2, 8).Value = "http://career.cmbchina.com/Campus/Position.aspx?id=10027"

# Row 3 - overseas management trainee, investment banking & financial markets
$ws.Cells.Item(3, 1).Value = "Not available"
Set-TextValue $ws.Cells.Item(3, 2) "2016-12-04"
$ws.Cells.Item(3, 3).Value = "深圳（轮岗地：深圳、北京、上海、广州、武汉、苏州、南京、天津、重庆）"
$ws.Cells.Item(3, 4).Value = "总行"
$ws.Cells.Item(3, 5).Value = "管理培训生（海外专场--投行与金融市场方向）"
$ws.Cells.Item(3, 6).Value = "C"
$ws.Cells.Item(3, 7).Value = "cmbchina"
$ws.Cells.Item(3, 8).Value = "http://career.cmbchina.com/Campus/Position.aspx?id=10295"

# Row 4 - overseas management trainee, retail finance
$ws.Cells.Item(4, 1).Value = "Not available"
Set-TextValue $ws.Cells.Item(4, 2) "2016-12-04"
$ws.Cells.Item(4, 3).Value = "深圳（轮岗地：深圳、北京、上海、广州、武汉、苏州、南京、天津、重庆）"
$ws.Cells.Item(4, 4).Value = "总行"
$ws.Cells.Item(4, 5).Value = "管理培训生（海外专场--零售金融方向）"
$ws.Cells.Item(4, 6).Value = "C"
$ws.Cells.Item(4, 7).Value = "cmbchina"
$ws.Cells.Item(4, 8).Value = "http://career.cmbchina.com/Campus/Position.aspx?id=10296"

# Row 5 - overseas management trainee, corporate finance
$ws.Cells.Item(5, 1).Value = "Not available"
Set-TextValue $ws.Cells.Item(5, 2) "2016-12-04"
$ws.Cells.Item(5, 3).Value = "深圳（轮岗地：深圳、北京、上海、广州、武汉、苏州、南京、天津、重庆）"
$ws.Cells.Item(5, 4).Value = "总行"
$ws.Cells.Item(5, 5).Value = "管理培训生（海外专场--公司金融方向）"
$ws.Cells.Item(5, 6).Value = "C"
$ws.Cells.Item(5, 7).Value = "cmbchina"
$ws.Cells.Item(5, 8).Value = "http://career.cmbchina.com/Campus/Position.aspx?id=10297"

# Keep the first sheet as the active one (matches original activeTab state)
$wb.Worksheets.Item(1).Activate()
